$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "sign up" / "create account" test data (row 21-23) ---
$ws.Range("A22").Value = "firstname"
$ws.Range("B22").Value = "lastname"
$ws.Range("C22").Value = "email"
$ws.Range("D22").Value = "newpassword"
$ws.Range("E22").Value = "confirmpassword"
$ws.Range("F22").Value = "phone"

$ws.Range("A23").Value = "Testing5"
$ws.Range("B23").Value = "Testing5"
$ws.Range("C23").Value = "Testing5@gmail.com"
$ws.Range("D23").Value = "Creating New Account"
$ws.Range("E23").Value = "Creating New Account"
$ws.Range("F23").Value = 8056037388

# --- Registered-user login data (row 26-28) ---
$ws.Range("A27").Value = "email"
$ws.Range("B27").Value = "password"

$ws.Range("A28").Value = "Testing4@gmail.com"
$ws.Range("B28").Value = "Creating New Account"

# --- Labels added after the table bodies were filled in ---
$ws.Range("A21").Value = "Newuser"
$ws.Range("A26").Value = "Registereduser"

# --- Hyperlinks for the two e-mail addresses ---
$ws.Hyperlinks.Add($ws.Range("C23"), "mailto:Testing5@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A28"), "mailto:Testing4@gmail.com") | Out-Null

# --- Column F needs to fit the phone number that was just entered ---
$ws.Columns("F").ColumnWidth = 10.14

# --- Selection / scroll position left by the editing session ---
$ws.Range("H29").Select()

Write-Host "done"
